$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.705.72"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "1.600.69"
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'211.39"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "'19.64"
$ws.Range("D11").Value = "'0.0845"
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("D12").Value = "1.825.44"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").Value = "1.600.46"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").Value = "'4.05"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("D16").Value = "'65.00"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "'210.07"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").Value = "'7.16"
$ws.Range("E20").Value = "  +2.12%  "
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("E22").Value = "  -3.17%  "
$ws.Range("D23").Value = "'8.97"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'143.64"
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("D25").Value = "'1.01"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("D28").Value = "'15.35"
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("E29").Value = "  -1.34%  "
$ws.Range("D30").Value = "'1.16"
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("D33").Value = "1.286.84"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  +0.73%  "
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("D36").Value = "'0.601"
$ws.Range("E36").Value = "  -3.28%  "
$ws.Range("E37").Value = "  +11.14%  "
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("E39").Value = "  -0.56%  "
$ws.Range("D40").Value = "'5.39"
$ws.Range("E40").Value = "  -2.04%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "'2.19"
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.784"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").Value = "'62.81"
$ws.Range("E43").Value = "  -1.24%  "
$ws.Range("D44").Value = "1.737.47"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").Value = "'90.48"
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("E46").Value = "  -1.32%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("E51").Value = "  +0.89%  "
